# Updates the cryptos list (prices and 1h volume % changes) to reflect the
# latest scrape. Two pairs of adjacent rows (11/12 and 29/30, plus the
# 13/14 and 31/32 pairs that follow them) swapped rank order, so their
# Coin/Link/Price/Volume values are exchanged accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.034.05'
$ws.Range("E2").Value = '  -6.28%  '
$ws.Range("D3").Value = '2.891.52'
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '545.87'
$ws.Range("E5").Value = '  -2.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '123.26'
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.500'
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("D9").Value = '2.887.20'
$ws.Range("E9").Value = '  -3.26%  '
$ws.Range("E10").Value = '  -9.82%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.64'
$ws.Range("E11").Value = '  -10.38%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.435'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.43'
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000208'
$ws.Range("E14").Value = '  -6.33%  '
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '3.359.90'
$ws.Range("E16").Value = '  -3.72%  '
$ws.Range("D17").Value = '2.880.58'
$ws.Range("E17").Value = '  -3.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.48'
$ws.Range("E18").Value = '  +5.53%  '
$ws.Range("D19").Value = '56.980.94'
$ws.Range("E19").Value = '  -6.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '404.83'
$ws.Range("E20").Value = '  -6.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.77'
$ws.Range("E21").Value = '  -2.44%  '
$ws.Range("E22").Value = '  +1.72%  '
$ws.Range("E23").Value = '  -4.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.60'
$ws.Range("E24").Value = '  -2.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '77.19'
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.44'
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.17'
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.91'
$ws.Range("E30").Value = '  +2.11%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '24.59'
$ws.Range("E31").Value = '  -2.86%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.92'
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0979'
$ws.Range("E33").Value = '  +5.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.909'
$ws.Range("E34").Value = '  -4.80%  '
$ws.Range("E35").Value = '  -3.44%  '
$ws.Range("E36").Value = '  -11.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '47.76'
$ws.Range("E37").Value = '  -4.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.26'
$ws.Range("E38").Value = '  +6.53%  '
$ws.Range("D39").Value = '0.0₃0620'
$ws.Range("E39").Value = '  -7.31%  '
$ws.Range("E40").Value = '  -1.63%  '
$ws.Range("E41").Value = '  -5.66%  '
$ws.Range("D42").Value = '2.623.40'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '358.48'
$ws.Range("E43").Value = '  -4.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.38'
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '119.59'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.227'
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.67'
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("E51").Value = '  -4.06%  '
